# Edit script: apply the changes described by the diff
$p = $ppt.ActivePresentation

# 1) Handout Master date placeholder: 27/06/2022 -> 28/06/2022
$hmDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$hmDate.Text = "28/06/2022"

# 2) Notes Master date placeholder: 27/06/2022 -> 28/06/2022
$nmDate = $p.NotesMaster.HeadersFooters.DateAndTime
$nmDate.Text = "28/06/2022"

# 3) Slide 12, table "Tableau 3", row 27 / col 2: clear placeholder text
#    "<OUI SI APE/ NON SINON>" -> empty (cell keeps its paragraph/run formatting)
$s12 = $p.Slides.Item(12)
$tbl = $s12.Shapes.Item(1).Table
$cell = $tbl.Cell(27, 2)
$cell.Shape.TextFrame.TextRange.Text = ""

# 4) Slide 9, shape "ZoneTexte 40": "<GCA>" -> "<CPR1>" within the long run
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(12)
$tr9 = $shp9.TextFrame.TextRange
$run = $tr9.Characters(151, 233)
$run.Text = "(<NSF> dans cet exemple). Le produit est automatiquement remboursé par anticipation. Il verse alors l’intégralité du capital initial majorée d’un <GC> de <CPN> par <F0> <F2> depuis le <DDCI>, soit un gain de <CPR1> dans notre exemple."
